$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing rows (13 -> 10 data rows total)
$ws.Rows("11:13").Delete()

# New data for rows 3 through 10 (rows 1 and 2 stay unchanged)
$data = @(
    @("In wich country is Barcelona located?", "Spain", "Location"),
    @("Where do most people speak italian?", "Italy", "Location"),
    @("Who was the F1 World Champion in 2022?", "Max Verstappen", "Person"),
    @("Who is the mayor of Innsbruck? ", "Georg Willi", "Person"),
    @("Who founded Facebook?", "Marc Zuckerberg", "Person"),
    @("When did Miachel Schumacher win his first F1 World Drivers Title?", 1994, "Year"),
    @("When did Miachel Schumacher win his 3rd F1 World Drivers Title?", 2000, "Year"),
    @("When has Chelsea last won the Champions League?", 2021, "Year")
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

$ws.Range("A1").Select()
$ws.Range("D20").Select()
